# Hotfix: fill in the remaining "Points" (column B) values for the
# checklist items that were already scored in column C but still had an
# empty Points cell, and move the active selection/scroll position to
# where the reviewer left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B24").Value = 3   # Edit user data, Profile-Page
$ws.Range("B26").Value = 3   # Trading
$ws.Range("B29").Value = 4   # Draw possible
$ws.Range("B36").Value = 3   # Mandatory Unique Feature
$ws.Range("B42").Value = 5   # Integration Tests

# B54 (Sum Points) recalculates automatically from these via its formula.

# Update the view: scroll so row 12 is the top visible row, and leave the
# selection on E42.
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E42").Select()
